$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Cotton Deck Mop Refill 8.1 oz."
$ws.Range("E10").Select()
